$wb = $excel.ActiveWorkbook

# --- Sheet "ev_charging_uc": update the two day/night timeslice order strings ---
$wsUc = $wb.Worksheets.Item("ev_charging_uc")
$wsUc.Range("C13").Value = "FaP,SaP,SaD,WaP,WaD,RaD,FaD,RaP"
$wsUc.Range("C14").Value = "WaP,SaN,WaN,FaP,SaP,RaN,FaN,RaP"

# --- Sheet "re_profiles": rotate the M/N (region, value) rows 4-7 up by one ---
$wsRe = $wb.Worksheets.Item("re_profiles")

$wsRe.Range("M4").Value = "W"
$wsRe.Range("N4").Value = 0.26654704733759038

$wsRe.Range("M5").Value = "F"
$wsRe.Range("N5").Value = 0.26274737453163755

$wsRe.Range("M6").Value = "S"
$wsRe.Range("N6").Value = 0.34545358594121062

$wsRe.Range("M7").Value = "R"
$wsRe.Range("N7").Value = 0.32525199218956147

$excel.Calculate()

$wb.Save()
